# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the
# a5368211-9e58-402a-b5ee-7f35c9e9b558 entry (row 4) on both the
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-28 03:58:20"
$wsZhCn.Range("G4").Value = "2016-01-28 03:59:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-28 03:58:31"
$wsDeDe.Range("G4").Value = "2016-01-28 03:59:17"
